# Update Deposit and Withdrawl Method Implementation
#
# 1. Swap the tab order of "amountDepositTest" and "amountwithdrawlTest"
#    (amountDepositTest now comes first).
# 2. Rewrite both sheets' content: common "TransactionAmount" header
#    (replacing the old "DepositAmount" / "WithDrawlAmount" headers) and
#    three transaction rows each.
# 3. Flip the runmode flag for the amountTransactionDetailsTest row on the
#    test_suite (TCID) sheet from Y to N, and update its selection.
# 4. Make "amountwithdrawlTest" the active/selected sheet (it used to be
#    "amountTransactionDetailsTest") -- done LAST so it sticks as the
#    workbook's active tab.

$wb = $excel.ActiveWorkbook

# --- 1. Reorder sheets: amountDepositTest before amountwithdrawlTest ---
# NOTE: worksheet object handles returned by Item() appear to be
# position-bound, so after Move() changes tab order we must re-resolve
# the sheets we want to edit by name again (stale handles now point at
# whatever sheet occupies their *original* slot/index).
$wsDeposit  = $wb.Worksheets.Item("amountDepositTest")
$wsWithdraw = $wb.Worksheets.Item("amountwithdrawlTest")
$wsDeposit.Move($wsWithdraw)

$wsDeposit  = $wb.Worksheets.Item("amountDepositTest")
$wsWithdraw = $wb.Worksheets.Item("amountwithdrawlTest")

# --- 2a. amountDepositTest: header + 3 rows (100, 200, 700) ---
$wsDeposit.Range("B1").Value = "TransactionAmount"

$wsDeposit.Range("A2").Value = "Dharmendra Pal"
$wsDeposit.Range("B2").Value = 100
$wsDeposit.Range("C2").Value = "Y"

$wsDeposit.Range("A3").Value = "Dharmendra Pal"
$wsDeposit.Range("B3").Value = 200
$wsDeposit.Range("C3").Value = "Y"

$wsDeposit.Range("A4").Value = "Dharmendra Pal"
$wsDeposit.Range("B4").Value = 700
$wsDeposit.Range("C4").Value = "Y"

$wsDeposit.Range("B9").Select() | Out-Null

# --- 2b. amountwithdrawlTest: header + 3 rows (200, 400, 300) ---
$wsWithdraw.Range("B1").Value = "TransactionAmount"

$wsWithdraw.Range("A2").Value = "Dharmendra Pal"
$wsWithdraw.Range("B2").Value = 200
$wsWithdraw.Range("C2").Value = "Y"

$wsWithdraw.Range("A3").Value = "Dharmendra Pal"
$wsWithdraw.Range("B3").Value = 400
$wsWithdraw.Range("C3").Value = "Y"

$wsWithdraw.Range("A4").Value = "Dharmendra Pal"
$wsWithdraw.Range("B4").Value = 300
$wsWithdraw.Range("C4").Value = "Y"

# --- 3. test_suite (TCID) sheet: runmode N for amountTransactionDetailsTest, selection B10 ---
$wsSuite = $wb.Worksheets.Item("test_suite")
$wsSuite.Range("B9").Value = "N"
$wsSuite.Range("B10").Select() | Out-Null

# --- 4. amountwithdrawlTest becomes the active sheet/selection (do this LAST) ---
$wsWithdraw = $wb.Worksheets.Item("amountwithdrawlTest")
$wsWithdraw.Activate()
$wsWithdraw.Range("B11").Select() | Out-Null
